# Update Hyperion_Profits market-price derived columns (H-N) across leve sheets.
# Generated from scheduled market data refresh.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 678.46155
$ws.Range("I2").Value = 282.3
$ws.Range("J2").Value = 1999
$ws.Range("K2").Value = 282.3
$ws.Range("L2").Value = 1999
$ws.Range("M2").Value = -169.3
$ws.Range("N2").Value = -2225
$ws.Range("H62").Value = 77840.07000000001
$ws.Range("I62").Value = 253062.25
$ws.Range("J62").Value = 7751.2
$ws.Range("K62").Value = 253062.25
$ws.Range("L62").Value = 7751.2
$ws.Range("M62").Value = -252438.25
$ws.Range("N62").Value = -8999.200000000001
$ws.Range("H64").Value = 12799.667
$ws.Range("J64").Value = 7221.6665
$ws.Range("L64").Value = 7221.6665
$ws.Range("N64").Value = -7717.6665
$ws.Range("H65").Value = 77840.07000000001
$ws.Range("I65").Value = 253062.25
$ws.Range("J65").Value = 7751.2
$ws.Range("K65").Value = 1265311.25
$ws.Range("L65").Value = 38756
$ws.Range("M65").Value = -1262191.25
$ws.Range("N65").Value = -44996
$ws.Range("H67").Value = 12799.667
$ws.Range("J67").Value = 7221.6665
$ws.Range("L67").Value = 7221.6665
$ws.Range("N67").Value = -8937.666499999999
$ws.Range("H86").Value = 5630.3477
$ws.Range("I86").Value = 5281.1665
$ws.Range("J86").Value = 5753.5884
$ws.Range("K86").Value = 5281.1665
$ws.Range("L86").Value = 5753.5884
$ws.Range("M86").Value = -4158.1665
$ws.Range("N86").Value = -7999.5884
$ws.Range("H89").Value = 5630.3477
$ws.Range("I89").Value = 5281.1665
$ws.Range("J89").Value = 5753.5884
$ws.Range("K89").Value = 26405.8325
$ws.Range("L89").Value = 28767.942
$ws.Range("M89").Value = -20789.8325
$ws.Range("N89").Value = -39999.942
$ws.Range("H106").Value = 20949.6
$ws.Range("I106").Value = 25312
$ws.Range("K106").Value = 25312
$ws.Range("M106").Value = -24681
$ws.Range("H112").Value = 6310.5654
$ws.Range("I112").Value = 2003.6666
$ws.Range("K112").Value = 6010.9998
$ws.Range("M112").Value = -4902.9998
$ws.Range("H132").Value = 22730416
$ws.Range("I132").Value = 27029850
$ws.Range("J132").Value = 4838.5713
$ws.Range("K132").Value = 81089550
$ws.Range("L132").Value = 14515.7139
$ws.Range("M132").Value = -81087020
$ws.Range("N132").Value = -19575.7139
$ws.Range("H135").Value = 637.8919
$ws.Range("J135").Value = 805.5
$ws.Range("L135").Value = 7249.5
$ws.Range("N135").Value = -12319.5
$ws.Range("H138").Value = 4215.523
$ws.Range("I138").Value = 2683.3125
$ws.Range("J138").Value = 5091.0713
$ws.Range("K138").Value = 8049.9375
$ws.Range("L138").Value = 15273.2139
$ws.Range("M138").Value = -2909.9375
$ws.Range("N138").Value = -25553.2139
$ws.Range("H141").Value = 4996.9287
$ws.Range("I141").Value = 5285.346
$ws.Range("K141").Value = 15856.038
$ws.Range("M141").Value = -10676.038

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3581.2615
$ws.Range("I32").Value = 2293.9075
$ws.Range("J32").Value = 9901
$ws.Range("K32").Value = 2293.9075
$ws.Range("L32").Value = 9901
$ws.Range("M32").Value = -2006.9075
$ws.Range("N32").Value = -10475
$ws.Range("H80").Value = 37999.5
$ws.Range("J80").Value = 50000
$ws.Range("L80").Value = 50000
$ws.Range("N80").Value = -51996
$ws.Range("H83").Value = 37999.5
$ws.Range("J83").Value = 50000
$ws.Range("L83").Value = 150000
$ws.Range("N83").Value = -159984

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 44992.5
$ws.Range("J132").Value = 44992.5
$ws.Range("L132").Value = 44992.5
$ws.Range("N132").Value = -55112.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 17039.863
$ws.Range("I31").Value = 8298.786
$ws.Range("K31").Value = 8298.786
$ws.Range("M31").Value = -8003.786
$ws.Range("H34").Value = 17039.863
$ws.Range("I34").Value = 8298.786
$ws.Range("K34").Value = 8298.786
$ws.Range("M34").Value = -8096.786
$ws.Range("H99").Value = 4883.385
$ws.Range("I99").Value = 5500
$ws.Range("K99").Value = 5500
$ws.Range("M99").Value = -4002
$ws.Range("H126").Value = 4883.385
$ws.Range("I126").Value = 5500
$ws.Range("K126").Value = 16500
$ws.Range("M126").Value = -14030

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 532
$ws.Range("I17").Value = 260.7143
$ws.Range("J17").Value = 1165
$ws.Range("K17").Value = 782.1428999999999
$ws.Range("L17").Value = 3495
$ws.Range("M17").Value = -613.1428999999999
$ws.Range("N17").Value = -3833
$ws.Range("H38").Value = 193.23077
$ws.Range("I38").Value = 169
$ws.Range("K38").Value = 507
$ws.Range("M38").Value = -160
$ws.Range("H68").Value = 2009.5385
$ws.Range("J68").Value = 2147.5715
$ws.Range("L68").Value = 6442.7145
$ws.Range("N68").Value = -8064.7145
$ws.Range("H71").Value = 2009.5385
$ws.Range("J71").Value = 2147.5715
$ws.Range("L71").Value = 19328.1435
$ws.Range("N71").Value = -27440.1435

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 16747
$ws.Range("J15").Value = 18996
$ws.Range("L15").Value = 18996
$ws.Range("N15").Value = -19572
$ws.Range("H81").Value = 16747
$ws.Range("J81").Value = 18996
$ws.Range("L81").Value = 18996
$ws.Range("N81").Value = -20992
$ws.Range("H84").Value = 16747
$ws.Range("J84").Value = 18996
$ws.Range("L84").Value = 56988
$ws.Range("N84").Value = -66972
$ws.Range("H97").Value = 994785.4399999999
$ws.Range("I97").Value = 1325130.5
$ws.Range("J97").Value = 3750
$ws.Range("K97").Value = 1325130.5
$ws.Range("L97").Value = 3750
$ws.Range("M97").Value = -1324634.5
$ws.Range("N97").Value = -4742
$ws.Range("H122").Value = 281959.25
$ws.Range("I122").Value = 374627.75
$ws.Range("J122").Value = 3953.75
$ws.Range("K122").Value = 1123883.25
$ws.Range("L122").Value = 11861.25
$ws.Range("M122").Value = -1121433.25
$ws.Range("N122").Value = -16761.25
$ws.Range("H126").Value = 3222974.5
$ws.Range("I126").Value = 1821012.6
$ws.Range("J126").Value = 5559577.5
$ws.Range("K126").Value = 5463037.800000001
$ws.Range("L126").Value = 16678732.5
$ws.Range("M126").Value = -5460567.800000001
$ws.Range("N126").Value = -16683672.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4727.9443
$ws.Range("I7").Value = 2259.25
$ws.Range("J7").Value = 9665.333000000001
$ws.Range("K7").Value = 2259.25
$ws.Range("L7").Value = 9665.333000000001
$ws.Range("M7").Value = -2147.25
$ws.Range("N7").Value = -9889.333000000001
$ws.Range("H122").Value = 6088.636
$ws.Range("J122").Value = 8157.8
$ws.Range("L122").Value = 24473.4
$ws.Range("N122").Value = -29373.4
$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("L123").ClearContents()
$ws.Range("M123").ClearContents()
$ws.Range("N123").Value = 0
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").ClearContents()
$ws.Range("N124").Value = 0
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").ClearContents()
$ws.Range("N125").Value = 0
$ws.Range("H126").Value = 4727.9443
$ws.Range("I126").Value = 2259.25
$ws.Range("J126").Value = 9665.333000000001
$ws.Range("K126").Value = 6777.75
$ws.Range("L126").Value = 28995.999
$ws.Range("M126").Value = -4307.75
$ws.Range("N126").Value = -33935.999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 209699.6
$ws.Range("J41").Value = 209699.6
$ws.Range("L41").Value = 209699.6
$ws.Range("N41").Value = -210479.6
$ws.Range("H45").Value = 11279.875
$ws.Range("J45").Value = 11310
$ws.Range("L45").Value = 11310
$ws.Range("N45").Value = -12292
$ws.Range("H54").Value = 30077
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("H81").Value = 55561944
$ws.Range("I81").Value = 166666670
$ws.Range("K81").Value = 333333340
$ws.Range("M81").Value = -333332279
$ws.Range("H84").Value = 55561944
$ws.Range("I84").Value = 166666670
$ws.Range("K84").Value = 1666666700
$ws.Range("M84").Value = -1666661396
$ws.Range("H100").Value = 1949.6154
$ws.Range("I100").Value = 1918.125
$ws.Range("K100").Value = 3836.25
$ws.Range("M100").Value = -3295.25
$ws.Range("H126").Value = 2390.375
$ws.Range("I126").Value = 2110.5
$ws.Range("J126").Value = 3789.75
$ws.Range("K126").Value = 6331.5
$ws.Range("L126").Value = 11369.25
$ws.Range("M126").Value = -3861.5
$ws.Range("N126").Value = -16309.25

Write-Host "Updated cells: 219 set, 5 cleared across 8 sheets."
